$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D:E").Insert()
$ws.Range("F:G").Copy()
$ws.Range("D:E").PasteSpecial(-4122)
$ws.Range("D5:E6").Clear()
$ws.Range("D36:E37").Clear()
$ws.Range("D78:E79").Clear()
$ws.Range("D7").Value2 = 43465
$ws.Range("E7").Value2 = 43373
$ws.Range("D8").Value2 = 141800
$ws.Range("E8").Value2 = 83900
$ws.Range("F8").Value2 = 154400
$ws.Range("G8").Value2 = 49100
$ws.Range("H8").Value2 = 58500
$ws.Range("I8").Value2 = 31000
$ws.Range("J8").Value2 = 40200
$ws.Range("D9").Value2 = 103500
$ws.Range("E9").Value2 = 58700
$ws.Range("F9").Value2 = 111400
$ws.Range("G9").Value2 = 34600
$ws.Range("H9").Value2 = 39500
$ws.Range("I9").Value2 = 21000
$ws.Range("J9").Value2 = 28300
$ws.Range("D10").Value2 = 38300
$ws.Range("E10").Value2 = 25200
$ws.Range("F10").Value2 = 43000
$ws.Range("G10").Value2 = 14500
$ws.Range("H10").Value2 = 19000
$ws.Range("I10").Value2 = 10000
$ws.Range("J10").Value2 = 11900
$ws.Range("D12").Value2 = 6200
$ws.Range("E12").Value2 = 5000
$ws.Range("F12").Value2 = 7300
$ws.Range("H12").Value2 = 3800
$ws.Range("J12").Value2 = 3300
$ws.Range("D13").Value2 = 0
$ws.Range("E13").Value2 = 0
$ws.Range("D14").Value2 = 0
$ws.Range("E14").Value2 = 0
$ws.Range("D15").Value2 = 0
$ws.Range("E15").Value2 = 0
$ws.Range("D17").Value2 = 132900
$ws.Range("E17").Value2 = 92700
$ws.Range("F17").Value2 = 142600
$ws.Range("G17").Value2 = 44300
$ws.Range("H17").Value2 = 50600
$ws.Range("I17").Value2 = 26400
$ws.Range("J17").Value2 = 37000
$ws.Range("D18").Value2 = 8900
$ws.Range("E18").Value2 = -8800
$ws.Range("F18").Value2 = 11800
$ws.Range("G18").Value2 = 4800
$ws.Range("H18").Value2 = 7900
$ws.Range("I18").Value2 = 4600
$ws.Range("J18").Value2 = 3200
$ws.Range("D20").Value2 = 500
$ws.Range("E20").Value2 = 400
$ws.Range("F20").Value2 = 400
$ws.Range("D21").Value2 = "NA"
$ws.Range("E21").Value2 = "NA"
$ws.Range("F21").Value2 = 12300
$ws.Range("H21").Value2 = 8200
$ws.Range("J21").Value2 = 3500
$ws.Range("D22").Value2 = 0
$ws.Range("E22").Value2 = 0
$ws.Range("D23").Value2 = 9500
$ws.Range("E23").Value2 = -8400
$ws.Range("F23").Value2 = 12200
$ws.Range("G23").Value2 = 5000
$ws.Range("H23").Value2 = 8000
$ws.Range("I23").Value2 = 4700
$ws.Range("J23").Value2 = 3300
$ws.Range("D24").Value2 = 1300
$ws.Range("E24").Value2 = 500
$ws.Range("F24").Value2 = 1700
$ws.Range("I24").Value2 = 700
$ws.Range("J24").Value2 = 500
$ws.Range("D25").Value2 = 0
$ws.Range("E25").Value2 = 0
$ws.Range("D26").Value2 = 8100
$ws.Range("E26").Value2 = -8900
$ws.Range("F26").Value2 = 10400
$ws.Range("G26").Value2 = 4400
$ws.Range("H26").Value2 = 7000
$ws.Range("I26").Value2 = 4100
$ws.Range("J26").Value2 = 2800
$ws.Range("D27").Value2 = 8100
$ws.Range("E27").Value2 = -9600
$ws.Range("F27").Value2 = 9000
$ws.Range("G27").Value2 = 4400
$ws.Range("H27").Value2 = 800
$ws.Range("I27").Value2 = 400
$ws.Range("J27").Value2 = 200
$ws.Range("D28").Value2 = 0
$ws.Range("E28").Value2 = 0
$ws.Range("D29").Value2 = 0
$ws.Range("E29").Value2 = 0
$ws.Range("D30").Value2 = 0
$ws.Range("E30").Value2 = 0
$ws.Range("D31").Value2 = 0
$ws.Range("E31").Value2 = 0
$ws.Range("D32").Value2 = -500
$ws.Range("E32").Value2 = -400
$ws.Range("F32").Value2 = -400
$ws.Range("D33").Value2 = 8100
$ws.Range("E33").Value2 = -9600
$ws.Range("F33").Value2 = 9000
$ws.Range("G33").Value2 = 4400
$ws.Range("H33").Value2 = 800
$ws.Range("I33").Value2 = 400
$ws.Range("J33").Value2 = 200
$ws.Range("D34").Value2 = 0
$ws.Range("E34").Value2 = 0
$ws.Range("D35").Value2 = 8100
$ws.Range("E35").Value2 = -9600
$ws.Range("F35").Value2 = 9000
$ws.Range("G35").Value2 = 4400
$ws.Range("H35").Value2 = 800
$ws.Range("I35").Value2 = 400
$ws.Range("J35").Value2 = 200
$ws.Range("D38").Value2 = 43465
$ws.Range("E38").Value2 = 43373
$ws.Range("D41").Value2 = 139500
$ws.Range("E41").Value2 = 138500
$ws.Range("F41").Value2 = 38100
$ws.Range("D42").Value2 = 25100
$ws.Range("E42").Value2 = "NA"
$ws.Range("F42").Value2 = "NA"
$ws.Range("G42").Value2 = "NA"
$ws.Range("H42").Value2 = "NA"
$ws.Range("I42").Value2 = "NA"
$ws.Range("J42").Value2 = "NA"
$ws.Range("D43").Value2 = 72000
$ws.Range("E43").Value2 = 57600
$ws.Range("F43").Value2 = 70400
$ws.Range("D44").Value2 = 34400
$ws.Range("E44").Value2 = 29800
$ws.Range("F44").Value2 = 24700
$ws.Range("D45").Value2 = 11300
$ws.Range("E45").Value2 = 13300
$ws.Range("F45").Value2 = 15700
$ws.Range("D46").Value2 = 282400
$ws.Range("E46").Value2 = 239200
$ws.Range("F46").Value2 = 148800
$ws.Range("D47").Value2 = 0
$ws.Range("E47").Value2 = 0
$ws.Range("D48").Value2 = 1700
$ws.Range("E48").Value2 = 500
$ws.Range("D49").Value2 = 0
$ws.Range("E49").Value2 = 0
$ws.Range("F49").Value2 = "NA"
$ws.Range("G49").Value2 = "NA"
$ws.Range("H49").Value2 = "NA"
$ws.Range("I49").Value2 = "NA"
$ws.Range("J49").Value2 = "NA"
$ws.Range("D50").Value2 = 0
$ws.Range("E50").Value2 = 0
$ws.Range("D51").Value2 = 0
$ws.Range("E51").Value2 = 0
$ws.Range("D52").Value2 = 1300
$ws.Range("E52").Value2 = 1500
$ws.Range("F52").Value2 = 1600
$ws.Range("D53").Value2 = 0
$ws.Range("E53").Value2 = 0
$ws.Range("D54").Value2 = 285400
$ws.Range("E54").Value2 = 241200
$ws.Range("F54").Value2 = 150900
$ws.Range("D57").Value2 = 81400
$ws.Range("E57").Value2 = 61600
$ws.Range("F57").Value2 = 73900
$ws.Range("D58").Value2 = 900
$ws.Range("E58").Value2 = 300
$ws.Range("D59").Value2 = 44100
$ws.Range("E59").Value2 = 29900
$ws.Range("F59").Value2 = 29100
$ws.Range("D60").Value2 = 126400
$ws.Range("E60").Value2 = 91800
$ws.Range("F60").Value2 = 103500
$ws.Range("D61").Value2 = 0
$ws.Range("E61").Value2 = 0
$ws.Range("D62").Value2 = 100
$ws.Range("E62").Value2 = 0
$ws.Range("D63").Value2 = 0
$ws.Range("E63").Value2 = 0
$ws.Range("D64").Value2 = 0
$ws.Range("E64").Value2 = 0
$ws.Range("D65").Value2 = 0
$ws.Range("E65").Value2 = 0
$ws.Range("D66").Value2 = 126900
$ws.Range("E66").Value2 = 92300
$ws.Range("F66").Value2 = 103500
$ws.Range("D68").Value2 = 0
$ws.Range("E68").Value2 = 0
$ws.Range("D69").Value2 = 0
$ws.Range("E69").Value2 = 0
$ws.Range("D70").Value2 = 0
$ws.Range("E70").Value2 = 0
$ws.Range("F70").Value2 = 23400
$ws.Range("D71").Value2 = 0
$ws.Range("E71").Value2 = 0
$ws.Range("D72").Value2 = -14200
$ws.Range("E72").Value2 = -22300
$ws.Range("F72").Value2 = -13400
$ws.Range("D73").Value2 = 0
$ws.Range("E73").Value2 = 0
$ws.Range("D74").Value2 = 0
$ws.Range("E74").Value2 = 0
$ws.Range("D75").Value2 = 0
$ws.Range("E75").Value2 = 0
$ws.Range("D76").Value2 = 158500
$ws.Range("E76").Value2 = 148900
$ws.Range("F76").Value2 = 24000
$ws.Range("D77").Value2 = 0
$ws.Range("E77").Value2 = 0
$ws.Range("D80").Value2 = 43465
$ws.Range("E80").Value2 = 43373
$ws.Range("D81").Value2 = 8100
$ws.Range("E81").Value2 = -9600
$ws.Range("F81").Value2 = 9000
$ws.Range("G81").Value2 = 4400
$ws.Range("H81").Value2 = 800
$ws.Range("I81").Value2 = 400
$ws.Range("J81").Value2 = 200
$ws.Range("D83").Value2 = "NA"
$ws.Range("E83").Value2 = "NA"
$ws.Range("F83").Value2 = 200
$ws.Range("G83").Value2 = "NA"
$ws.Range("H83").Value2 = "NA"
$ws.Range("I83").Value2 = "NA"
$ws.Range("J83").Value2 = "NA"
$ws.Range("D84").Value2 = 0
$ws.Range("E84").Value2 = 0
$ws.Range("D85").Value2 = 0
$ws.Range("E85").Value2 = 0
$ws.Range("D86").Value2 = 0
$ws.Range("E86").Value2 = 0
$ws.Range("D87").Value2 = 0
$ws.Range("E87").Value2 = 0
$ws.Range("D88").Value2 = 0
$ws.Range("E88").Value2 = 0
$ws.Range("D89").Value2 = "NA"
$ws.Range("E89").Value2 = "NA"
$ws.Range("F89").Value2 = -2700
$ws.Range("G89").Value2 = "NA"
$ws.Range("H89").Value2 = "NA"
$ws.Range("I89").Value2 = "NA"
$ws.Range("J89").Value2 = "NA"
$ws.Range("D91").Value2 = "NA"
$ws.Range("E91").Value2 = "NA"
$ws.Range("F91").Value2 = -900
$ws.Range("G91").Value2 = "NA"
$ws.Range("H91").Value2 = "NA"
$ws.Range("I91").Value2 = "NA"
$ws.Range("J91").Value2 = "NA"
$ws.Range("D92").Value2 = 0
$ws.Range("E92").Value2 = 0
$ws.Range("D93").Value2 = 0
$ws.Range("E93").Value2 = 0
$ws.Range("D94").Value2 = "NA"
$ws.Range("E94").Value2 = "NA"
$ws.Range("F94").Value2 = 3800
$ws.Range("G94").Value2 = "NA"
$ws.Range("H94").Value2 = "NA"
$ws.Range("I94").Value2 = "NA"
$ws.Range("J94").Value2 = "NA"
$ws.Range("D96").Value2 = 0
$ws.Range("E96").Value2 = 0
$ws.Range("D97").Value2 = 0
$ws.Range("E97").Value2 = 0
$ws.Range("D98").Value2 = 0
$ws.Range("E98").Value2 = 0
$ws.Range("D99").Value2 = 0
$ws.Range("E99").Value2 = 0
$ws.Range("D100").Value2 = "NA"
$ws.Range("E100").Value2 = "NA"
$ws.Range("F100").Value2 = -5100
$ws.Range("G100").Value2 = "NA"
$ws.Range("H100").Value2 = "NA"
$ws.Range("I100").Value2 = "NA"
$ws.Range("J100").Value2 = "NA"
$ws.Range("D101").Value2 = "NA"
$ws.Range("E101").Value2 = "NA"
$ws.Range("F101").Value2 = 600
$ws.Range("G101").Value2 = "NA"
$ws.Range("H101").Value2 = "NA"
$ws.Range("I101").Value2 = "NA"
$ws.Range("J101").Value2 = "NA"
$ws.Range("D102").Value2 = "NA"
$ws.Range("E102").Value2 = "NA"
$ws.Range("F102").Value2 = -3400
$ws.Range("G102").Value2 = "NA"
$ws.Range("H102").Value2 = "NA"
$ws.Range("I102").Value2 = "NA"
$ws.Range("J102").Value2 = "NA"
